$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values are plain text (non-numeric-looking), set directly
$ws.Range('D2').Value = '30.244.40'
$ws.Range('E2').Value = '  +3.30%  '
$ws.Range('D3').Value = '1.908.67'
$ws.Range('E3').Value = '  +0.34%  '
$ws.Range('E4').Value = '  -0.35%  '
$ws.Range('E5').Value = '  +3.72%  '
$ws.Range('E6').Value = '  -0.32%  '
$ws.Range('E7').Value = '  +0.59%  '
$ws.Range('E8').Value = '  +1.98%  '
$ws.Range('E9').Value = '  +0.26%  '
$ws.Range('E10').Value = '  +0.44%  '
$ws.Range('E11').Value = '  +0.10%  '
$ws.Range('E12').Value = '  +13.52%  '
$ws.Range('E13').Value = '  +3.27%  '
$ws.Range('D14').Value = '1.906.09'
$ws.Range('E14').Value = '  +0.19%  '
$ws.Range('E15').Value = '  +0.08%  '
$ws.Range('E16').Value = '  -0.35%  '
$ws.Range('E18').Value = '  +0.59%  '
$ws.Range('E20').Value = '  +2.81%  '
$ws.Range('E21').Value = '  -0.27%  '
$ws.Range('E22').Value = '  -0.70%  '
$ws.Range('D23').Value = '30.249.54'
$ws.Range('E23').Value = '  +3.27%  '
$ws.Range('E24').Value = '  +1.01%  '
$ws.Range('E25').Value = '  -0.32%  '
$ws.Range('D26').Value = '2.126.18'
$ws.Range('E26').Value = '  +0.29%  '
$ws.Range('E27').Value = '  +3.52%  '
$ws.Range('E28').Value = '  +0.60%  '
$ws.Range('E29').Value = '  -1.71%  '
$ws.Range('E30').Value = '  +2.15%  '
$ws.Range('E31').Value = '  +3.59%  '
$ws.Range('E32').Value = '  +1.16%  '
$ws.Range('E33').Value = '  +1.31%  '
$ws.Range('E34').Value = '  +2.73%  '
$ws.Range('E35').Value = '  +0.79%  '
$ws.Range('E36').Value = '  -0.69%  '
$ws.Range('E37').Value = '  +0.77%  '
$ws.Range('E38').Value = '  +2.03%  '
$ws.Range('E39').Value = '  -0.38%  '
$ws.Range('E40').Value = '  +5.78%  '
$ws.Range('E41').Value = '  -3.61%  '
$ws.Range('E42').Value = '  +0.48%  '
$ws.Range('E43').Value = '  -0.45%  '
$ws.Range('E44').Value = '  +1.30%  '
$ws.Range('E45').Value = '  +1.02%  '
$ws.Range('E46').Value = '  +1.06%  '
$ws.Range('E47').Value = '  +0.61%  '
$ws.Range('E48').Value = '  +1.19%  '
$ws.Range('E49').Value = '  +1.38%  '
$ws.Range('E50').Value = '  -1.87%  '
$ws.Range('E51').Value = '  +1.83%  '

# Price cells whose new values look like plain numbers: force them to remain
# text (matching the original inlineStr cell type) by applying a text number
# format before assignment, then resetting the cell style back to Normal so
# no stray formatting is left behind.
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '327.09'
$c.Style = 'Normal'
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '0.9999'
$c.Style = 'Normal'
$c = $ws.Range('D7')
$c.NumberFormat = '@'
$c.Value = '0.5153'
$c.Style = 'Normal'
$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '0.4019'
$c.Style = 'Normal'
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '0.08478'
$c.Style = 'Normal'
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '42.67'
$c.Style = 'Normal'
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '23.43'
$c.Style = 'Normal'
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '6.456'
$c.Style = 'Normal'
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '7.371'
$c.Style = 'Normal'
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '1.000'
$c.Style = 'Normal'
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '95.02'
$c.Style = 'Normal'
$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '0.00001113'
$c.Style = 'Normal'
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '0.06668'
$c.Style = 'Normal'
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '18.42'
$c.Style = 'Normal'
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '0.9999'
$c.Style = 'Normal'
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '5.998'
$c.Style = 'Normal'
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '11.28'
$c.Style = 'Normal'
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '2.209'
$c.Style = 'Normal'
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '21.69'
$c.Style = 'Normal'
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '161.30'
$c.Style = 'Normal'
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '2.397'
$c.Style = 'Normal'
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '129.68'
$c.Style = 'Normal'
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '1.099'
$c.Style = 'Normal'
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '3.749'
$c.Style = 'Normal'
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '0.02497'
$c.Style = 'Normal'
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '0.06570'
$c.Style = 'Normal'
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '0.2213'
$c.Style = 'Normal'
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '5.219'
$c.Style = 'Normal'
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '1.230'
$c.Style = 'Normal'
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '11.93'
$c.Style = 'Normal'
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '8.794'
$c.Style = 'Normal'
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '0.6516'
$c.Style = 'Normal'
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '0.6130'
$c.Style = 'Normal'
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '13.29'
$c.Style = 'Normal'
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '3.719'
$c.Style = 'Normal'
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '2.064'
$c.Style = 'Normal'
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '1.246'
$c.Style = 'Normal'
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '124.84'
$c.Style = 'Normal'
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '1.157'
$c.Style = 'Normal'
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '79.34'
$c.Style = 'Normal'
